$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.233906149864197
$ws.Range("B1").Value = 3.183047533035278
$ws.Range("C1").Value = 6.304697036743164
$ws.Range("D1").Value = 1.816049814224243
$ws.Range("E1").Value = 1.068257331848145
